$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BXS")

# Income Statement section - Earnings Before Interest And Taxes (row 21), 2011 column (J) -> NA
$ws.Range("J21").Value = "NA"

# Balance Sheet section
# Cash And Cash Equivalents (row 41), 2017 column (D)
$ws.Range("D41").Value = 154800
# Short Term Investments (row 42), 2017 column (D)
$ws.Range("D42").Value = 114500
# Property Plant and Equipment (row 48), 2017 column (D)
$ws.Range("D48").Value = 628200
# Goodwill (row 49), 2017 column (D)
$ws.Range("D49").Value = 336900
# Total Assets (row 54), 2017 column (D)
$ws.Range("D54").Value = 14809500
# Accounts Payable (row 57), 2017 column (D)
$ws.Range("D57").Value = 4700
# Total Liabilities (row 66), 2017 column (D)
$ws.Range("D66").Value = 13108300
# Retained Earnings (row 72), 2017 column (D)
$ws.Range("D72").Value = 1354000
# Total Stockholder Equity (row 76), 2017 column (D)
$ws.Range("D76").Value = 1701200

# Cash Flow Statement section - 2011 column (J) -> NA
# Depreciation (row 83)
$ws.Range("J83").Value = "NA"
# Total Cash Flows From Investing Activities (row 94)
$ws.Range("J94").Value = "NA"
# Total Cash Flows From Financing Activities (row 100)
$ws.Range("J100").Value = "NA"
# Effect Of Exchange Rate Changes (row 101)
$ws.Range("J101").Value = "NA"
